$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 912
$ws1.Range("F3").Value = 504
$ws1.Range("F4").Value = 504
$ws1.Range("F5").Value = 782
$ws1.Range("F6").Value = 185
$ws1.Range("F7").Value = 1378
$ws1.Range("F8").Value = 804
$ws1.Range("F9").Value = 456
$ws1.Range("F10").Value = 608
$ws1.Range("F11").Value = 160
$ws1.Range("F14").Value = 205
$ws1.Range("F15").Value = 104
$ws1.Range("F16").Value = 1526
$ws1.Range("F17").Value = 184
$ws1.Range("F18").Value = 25
$ws1.Range("F19").Value = 456
$ws1.Range("F20").Value = 51
$ws1.Range("F22").Value = 104
$ws1.Range("F23").Value = 626
$ws1.Range("F24").Value = 15
$ws1.Range("F25").Value = 204
$ws1.Range("F26").Value = 708
$ws1.Range("F28").Value = 1358
$ws1.Range("F29").Value = 114

# Sheet "演出" (Performance) - sheet2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 6

# Sheet "本地生活" (Local Life) - sheet3
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 387

# Sheet "全部类型" (All Types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 387
$ws4.Range("F3").Value = 912
$ws4.Range("F4").Value = 504
$ws4.Range("F5").Value = 504
$ws4.Range("F6").Value = 782
$ws4.Range("F7").Value = 185
$ws4.Range("F8").Value = 1378
$ws4.Range("F9").Value = 804
$ws4.Range("F12").Value = 456
$ws4.Range("F13").Value = 608
$ws4.Range("F15").Value = 160
$ws4.Range("F18").Value = 205
$ws4.Range("F19").Value = 104
$ws4.Range("F20").Value = 1526
$ws4.Range("F22").Value = 184
$ws4.Range("F23").Value = 25
$ws4.Range("F24").Value = 456
$ws4.Range("F25").Value = 51
$ws4.Range("F27").Value = 6
$ws4.Range("F28").Value = 104
$ws4.Range("F31").Value = 626
$ws4.Range("F36").Value = 15
$ws4.Range("F37").Value = 204
$ws4.Range("F38").Value = 708
$ws4.Range("F40").Value = 1358
$ws4.Range("F41").Value = 114
